# Gnai2-Cnr1.xlsx: refresh with new TPM-based NATMI numbers.
# Rows 2-9 get updated values (ligand/receptor cluster pairing + recomputed
# specificity/weight metrics); 4 new rows (10-13) are appended to complete
# the full Sending x Target cluster combination matrix (ECs/FAPs/MuSCs x
# ECs/FAPs/MuSCs, plus Resolving-Mac as a sending cluster).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Cnr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1408123333333333
$ws.Range("N2").Value = 0.422437
$ws.Range("O2").Value = 0.01914604775478091
$ws.Range("P2").Value = 0.01914604775478091
$ws.Range("Q2").Value = 23.92622013203634
$ws.Range("R2").Value = 215.335981188327
$ws.Range("S2").Value = 0.008503113879080232
$ws.Range("T2").Value = 0.008503113879080234

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Cnr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.423734333333333
$ws.Range("N3").Value = 7.271203
$ws.Range("O3").Value = 0.3295516251481434
$ws.Range("P3").Value = 0.3295516251481434
$ws.Range("Q3").Value = 411.8304116417904
$ws.Range("R3").Value = 3706.473704776113
$ws.Range("S3").Value = 0.1463599711836554
$ws.Range("T3").Value = 0.1463599711836554

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Cnr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.790095666666667
$ws.Range("N4").Value = 14.370287
$ws.Range("O4").Value = 0.6513023270970757
$ws.Range("P4").Value = 0.6513023270970758
$ws.Range("Q4").Value = 813.9122522945197
$ws.Range("R4").Value = 7325.210270650678
$ws.Range("S4").Value = 0.2892554081107153
$ws.Range("T4").Value = 0.2892554081107154

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Cnr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 68.382243
$ws.Range("H5").Value = 205.146729
$ws.Range("I5").Value = 0.1787346690539575
$ws.Range("J5").Value = 0.1787346690539575
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1408123333333333
$ws.Range("N5").Value = 0.422437
$ws.Range("O5").Value = 0.01914604775478091
$ws.Range("P5").Value = 0.01914604775478091
$ws.Range("Q5").Value = 9.629063195397
$ws.Range("R5").Value = 86.661568758573
$ws.Range("S5").Value = 0.003422062509142033
$ws.Range("T5").Value = 0.003422062509142033

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Cnr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 68.382243
$ws.Range("H6").Value = 205.146729
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.423734333333333
$ws.Range("N6").Value = 7.271203
$ws.Range("O6").Value = 0.3295516251481434
$ws.Range("P6").Value = 0.3295516251481434
$ws.Range("Q6").Value = 165.740390149443
$ws.Range("R6").Value = 1491.663511344987
$ws.Range("S6").Value = 0.05890230065704728
$ws.Range("T6").Value = 0.05890230065704728

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Cnr1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 68.382243
$ws.Range("H7").Value = 205.146729
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.790095666666667
$ws.Range("N7").Value = 14.370287
$ws.Range("O7").Value = 0.6513023270970757
$ws.Range("P7").Value = 0.6513023270970758
$ws.Range("Q7").Value = 327.557485871247
$ws.Range("R7").Value = 2948.017372841223
$ws.Range("S7").Value = 0.1164103058877682
$ws.Range("T7").Value = 0.1164103058877682

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Cnr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 53.27463399999999
$ws.Range("H8").Value = 159.823902
$ws.Range("I8").Value = 0.1392470275793777
$ws.Range("J8").Value = 0.1392470275793778
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1408123333333333
$ws.Range("N8").Value = 0.422437
$ws.Range("O8").Value = 0.01914604775478091
$ws.Range("P8").Value = 0.01914604775478091
$ws.Range("Q8").Value = 7.501725521019333
$ws.Range("R8").Value = 67.51552968917399
$ws.Range("S8").Value = 0.002666030239746061
$ws.Range("T8").Value = 0.002666030239746062

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Cnr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 53.27463399999999
$ws.Range("H9").Value = 159.823902
$ws.Range("I9").Value = 0.1392470275793777
$ws.Range("J9").Value = 0.1392470275793778
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.423734333333333
$ws.Range("N9").Value = 7.271203
$ws.Range("O9").Value = 0.3295516251481434
$ws.Range("P9").Value = 0.3295516251481434
$ws.Range("Q9").Value = 129.1235595215673
$ws.Range("R9").Value = 1162.112035694106
$ws.Range("S9").Value = 0.04588908423583228
$ws.Range("T9").Value = 0.04588908423583229

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Cnr1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.790095666666667
$ws.Range("N10").Value = 14.370287
$ws.Range("O10").Value = 0.6513023270970757
$ws.Range("P10").Value = 0.6513023270970758
$ws.Range("Q10").Value = 255.1905934666526
$ws.Range("R10").Value = 2296.715341199874
$ws.Range("S10").Value = 0.0906919131037994
$ws.Range("T10").Value = 0.09069191310379943

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Gnai2"
$ws.Range("C11").Value = "Cnr1"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 91.01828266666666
$ws.Range("H11").Value = 273.054848
$ws.Range("I11").Value = 0.2378998101932138
$ws.Range("J11").Value = 0.2378998101932138
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1408123333333333
$ws.Range("N11").Value = 0.422437
$ws.Range("O11").Value = 0.01914604775478091
$ws.Range("P11").Value = 0.01914604775478091
$ws.Range("Q11").Value = 12.81649675828622
$ws.Range("R11").Value = 115.348470824576
$ws.Range("S11").Value = 0.004554841126812587
$ws.Range("T11").Value = 0.004554841126812587

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Gnai2"
$ws.Range("C12").Value = "Cnr1"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 91.01828266666666
$ws.Range("H12").Value = 273.054848
$ws.Range("I12").Value = 0.2378998101932138
$ws.Range("J12").Value = 0.2378998101932138
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.423734333333333
$ws.Range("N12").Value = 7.271203
$ws.Range("O12").Value = 0.3295516251481434
$ws.Range("P12").Value = 0.3295516251481434
$ws.Range("Q12").Value = 220.6041366602382
$ws.Range("R12").Value = 1985.437229942144
$ws.Range("S12").Value = 0.07840026907160846
$ws.Range("T12").Value = 0.07840026907160846

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Gnai2"
$ws.Range("C13").Value = "Cnr1"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 91.01828266666666
$ws.Range("H13").Value = 273.054848
$ws.Range("I13").Value = 0.2378998101932138
$ws.Range("J13").Value = 0.2378998101932138
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.790095666666667
$ws.Range("N13").Value = 14.370287
$ws.Range("O13").Value = 0.6513023270970757
$ws.Range("P13").Value = 0.6513023270970758
$ws.Range("Q13").Value = 435.9862813890418
$ws.Range("R13").Value = 3923.876532501376
$ws.Range("S13").Value = 0.1549446999947928
$ws.Range("T13").Value = 0.1549446999947928
